# [PHOENIX-5917] updated fee details in legacy Trade License
#
# This script:
#  1. Renames the "legacyDetails" sheet to "legencyDetails" (typo fix, matches workbook.xml diff).
#  2. Rewrites the "Fee Details" header row of that sheet with the new
#     amount1 / amountB / amountC / amountD / amountE / amountF headers
#     (columns B:G), reusing the existing "dataName" header text in A1.
#  3. Adds a second data row (legency Trade / 10 / 20 / 30 / 40 / 50 / 60)
#     spanning A2:G2.
#  4. Applies a Text ("@") number format to column B and to the full new
#     data row, matching the extra cellXfs entry added upstream.
#  5. Re-points the selection/active cell on both the "tradeDetails" sheet
#     (B2 -> F6) and the legacy sheet (B11 -> E12), and makes sure the
#     legacy sheet stays the active tab/sheet, same as before the edit.

$wb = $excel.ActiveWorkbook

# --- 1. Rename the legacy sheet -------------------------------------------------
$ws = $wb.Worksheets.Item("legacyDetails")
$ws.Name = "legencyDetails"

# --- 2. Header row (row 1) -------------------------------------------------------
$ws.Range("A1").Value = "dataName"
$ws.Range("B1").Value = "amount1"
$ws.Range("C1").Value = "amountB"
$ws.Range("D1").Value = "amountC"
$ws.Range("E1").Value = "amountD"
$ws.Range("F1").Value = "amountE"
$ws.Range("G1").Value = "amountF"

# --- 3. Data row (row 2) ----------------------------------------------------------
$ws.Range("A2").Value = "legency Trade"
$ws.Range("B2").Value = 10
$ws.Range("C2").Value = 20
$ws.Range("D2").Value = 30
$ws.Range("E2").Value = 40
$ws.Range("F2").Value = 50
$ws.Range("G2").Value = 60

# --- 4. Text formatting for column B and the whole data row ----------------------
$ws.Columns.Item(2).NumberFormat = "@"
$ws.Range("A2:G2").NumberFormat = "@"

# --- 5. Selections -----------------------------------------------------------------
$wsTrade = $wb.Worksheets.Item("tradeDetails")
[void]$wsTrade.Range("F6").Select()

[void]$ws.Range("E12").Select()
